$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C7: was stored as text "12" -> now a numeric 12
$ws.Range("C7").Value = 12

# E7: was stored as text "123123" -> now a numeric 123123
$ws.Range("E7").Value = 123123

# F7: new numeric cell with value 20
$ws.Range("F7").Value = 20

# G7: was stored as text "22" -> now a numeric 999
$ws.Range("G7").Value = 999

# H7: text date value updated from "2023-09-07" to "2023-10-03",
# staying as text (not auto-converted to a date serial number).
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "2023-10-03"
$ws.Range("H7").Style = "Normal"
